# "Base del proyecto funcional"
#
# The workbook ("errores") is a small data-validation export with 3 data
# rows (2-4). This edit:
#   1. Corrects several fields on row 2 that the user re-typed
#      (regimen, tipo de identificacion, direccion, codigo municipio).
#   2. Replaces a batch of invalid/placeholder date values (which were
#      stored as out-of-range negative date serials, e.g. -36888) with
#      the literal text "1799-01-01" used by the validator as an
#      "invalid date" sentinel, across rows 2-4.
#   3. Flags two more fields on row 2 (45.2 / 45.4 "al inicio") with the
#      same yellow highlight previously (wrongly) sitting on "30. TFG
#      Inicial", which itself goes back to a normal, unhighlighted 300.
#   4. Fixes a couple of placeholder numeric codes that should have been
#      the habilitacion code text instead of the sentinel 999999.
#   5. Extends the custom date format used elsewhere in the sheet to
#      also carry a time portion.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Helper: write a value that LOOKS numeric/date-like as real, literal
# text (so it round-trips as a shared string, not an auto-parsed number
# or date serial). We stage the cell as Text, assign, then drop back to
# an un-styled ("Normal") cell so no residual formatting is left behind.
# ---------------------------------------------------------------------
function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# -----------------------------------------------------------------
# 1. Row 2 field corrections (actual re-typed data)
# -----------------------------------------------------------------
$ws.Range("B2").Value = "S"
$ws.Range("H2").Value = "MS"
$ws.Range("M2").Value = "CL 13 CRA 18   17-18 B FLORIDA MAGANGUE"
Set-TextValue $ws.Range("O2") "25280"

# Row 3's address cell was left pointing at the same (now corrected)
# text as row 2.
$ws.Range("M3").Value = "CL 13 CRA 18   17-18 B FLORIDA MAGANGUE"
# Row 4 got its own, slightly different, manual correction.
$ws.Range("M4").Value = "CL 13 CRA 18  17-18  B FLORIDA MAGANGUE"

# -----------------------------------------------------------------
# 2. Invalid-date sentinel cleanup: negative/garbage date serials
#    become the literal text "1799-01-01" across rows 2-4.
# -----------------------------------------------------------------
$dateCols = @("Q", "S", "AN", "AQ", "AX", "CP", "CT", "EL", "EO", "EQ")
foreach ($row in 2..4) {
    foreach ($col in $dateCols) {
        Set-TextValue $ws.Range("$col$row") "1799-01-01"
    }
}

# -----------------------------------------------------------------
# 3. Move the "needs review" yellow highlight from AD2 (30. TFG
#    Inicial, now a normal value of 300) onto AT2/AV2, which become
#    the text sentinel "m".
# -----------------------------------------------------------------
$ws.Range("AD2").Copy($ws.Range("AT2"))
$ws.Range("AD2").Copy($ws.Range("AV2"))
$ws.Range("AT2").Value = "m"
$ws.Range("AV2").Value = "m"

$ws.Range("AD2").Style = "Normal"
$ws.Range("AD2").Value = 300

# -----------------------------------------------------------------
# 4. EAPB/IPS habilitacion code sentinel fix (999999 -> real code
#    text) on rows 2 and 3.
# -----------------------------------------------------------------
Set-TextValue $ws.Range("EP2") "050010412701"
Set-TextValue $ws.Range("EP3") "050010412701"

# -----------------------------------------------------------------
# 5. Extend the custom date format with a time component.
# -----------------------------------------------------------------
$dateFormatCells = @("J2","P2","R2","AS2","DC2","J3","P3","R3","AS3","DC3","J4","P4","R4","AS4","DC4")
foreach ($ref in $dateFormatCells) {
    $ws.Range($ref).NumberFormat = "yyyy-mm-dd h:mm:ss"
}
